$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.031.72"
$ws.Range("E2").Value = "'  -0.14%  "
$ws.Range("D3").Value = "'2.056.37"
$ws.Range("E3").Value = "'  +0.27%  "
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'246.57"
$ws.Range("E5").Value = "'  -1.15%  "
$ws.Range("D6").Value = "'0.659"
$ws.Range("E6").Value = "'  -1.38%  "
$ws.Range("D7").Value = "'59.09"
$ws.Range("E7").Value = "'  -0.27%  "
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E9").Value = "'  -2.69%  "
$ws.Range("D10").Value = "'0.0783"
$ws.Range("E10").Value = "'  -1.66%  "
$ws.Range("E11").Value = "'  +2.44%  "
$ws.Range("D12").Value = "'15.46"
$ws.Range("E12").Value = "'  -3.61%  "
$ws.Range("D13").Value = "'0.889"
$ws.Range("E13").Value = "'  +6.11%  "
$ws.Range("D14").Value = "'2.357.66"
$ws.Range("E14").Value = "'  +0.39%  "
$ws.Range("D15").Value = "'5.72"
$ws.Range("E15").Value = "'  -0.13%  "
$ws.Range("D16").Value = "'2.055.09"
$ws.Range("E16").Value = "'  +0.19%  "
$ws.Range("D17").Value = "'18.24"
$ws.Range("E17").Value = "'  -0.85%  "
$ws.Range("D18").Value = "'37.003.19"
$ws.Range("E18").Value = "'  -0.08%  "
$ws.Range("D19").Value = "'73.96"
$ws.Range("E19").Value = "'  -1.92%  "
$ws.Range("E20").Value = "'  -0.96%  "
$ws.Range("E21").Value = "'  +0.65%  "
$ws.Range("D22").Value = "'238.60"
$ws.Range("E22").Value = "'  +0.24%  "
$ws.Range("E23").Value = "'  -0.03%  "
$ws.Range("E24").Value = "'  +1.54%  "
$ws.Range("D25").Value = "'10.23"
$ws.Range("E25").Value = "'  +7.83%  "
$ws.Range("D26").Value = "'170.51"
$ws.Range("E26").Value = "'  +0.69%  "
$ws.Range("E27").Value = "'  -2.37%  "
$ws.Range("E28").Value = "'  +0.14%  "
$ws.Range("D29").Value = "'5.51"
$ws.Range("E29").Value = "'  +14.91%  "
$ws.Range("E30").Value = "'  -0.92%  "
$ws.Range("E31").Value = "'  -0.69%  "
$ws.Range("D32").Value = "'4.72"
$ws.Range("E32").Value = "'  +4.36%  "
$ws.Range("E33").Value = "'  -1.06%  "
$ws.Range("D34").Value = "'2.38"
$ws.Range("E34").Value = "'  +6.57%  "
$ws.Range("E35").Value = "'  +0.03%  "
$ws.Range("E36").Value = "'  +5.58%  "
$ws.Range("D37").Value = "'0.0842"
$ws.Range("E37").Value = "'  -5.53%  "
$ws.Range("E38").Value = "'  +0.11%  "
$ws.Range("E39").Value = "'  +2.43%  "
$ws.Range("D40").Value = "'3.06"
$ws.Range("E40").Value = "'  -1.17%  "
$ws.Range("E41").Value = "'  +0.54%  "
$ws.Range("D42").Value = "'1.16"
$ws.Range("E42").Value = "'  +1.98%  "
$ws.Range("D43").Value = "'0.0959"
$ws.Range("E43").Value = "'  -10.78%  "
$ws.Range("D44").Value = "'97.89"
$ws.Range("E44").Value = "'  +0.86%  "
$ws.Range("D45").Value = "'17.02"
$ws.Range("E45").Value = "'  -3.93%  "
$ws.Range("D46").Value = "'1.304.69"
$ws.Range("E46").Value = "'  +1.23%  "
$ws.Range("E47").Value = "'  -6.06%  "
$ws.Range("E48").Value = "'  -0.38%  "
$ws.Range("D49").Value = "'6.82"
$ws.Range("E49").Value = "'  +0.21%  "
$ws.Range("D50").Value = "'2.246.09"
$ws.Range("E50").Value = "'  +0.96%  "
$ws.Range("D51").Value = "'44.80"
$ws.Range("E51").Value = "'  +3.26%  "

Write-Host "Updated cryptos list"